{"js": "// Insert two new bulleted list items after the \"Confirming the performance\n// of the algorithm tracking SOC and the remaining energy in the battery.\"\n// list entry, matching the same list paragraph style/numbering.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the last existing list item in this list (\"...remaining energy in\n// the battery.\") so the two new bullets are appended right after it.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"remaining energy in the battery\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find the anchor list paragraph.\");\n}\n\n// insertParagraph(\"After\") clones the anchor paragraph's formatting\n// (pStyle \"Listeavsnitt\", numPr ilvl 0 / numId 1, rPr lang en-GB), so the\n// new paragraphs automatically become list items in the same list.\nconst first = anchor.insertParagraph(\n  \"(Confirm the performance of the current sensor at low currents.)\",\n  \"After\"\n);\nawait context.sync();\n\nfirst.insertParagraph(\"Testing the new constant voltage state\", \"After\");\nawait context.sync();\n", "ps1": "# Insert two new bulleted list items after the \"Confirming the performance\n# of the algorithm tracking SOC and the remaining energy in the battery.\"\n# list entry, matching the same list paragraph style/numbering.\n\n$d = $word.ActiveDocument\n\n# Locate the last existing list item in this list (\"...remaining energy in\n# the battery.\") so the two new bullets are appended right after it.\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*remaining energy in the battery*\") {\n        $anchor = $p\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find the anchor list paragraph.\"\n}\n\n# InsertParagraphAfter() clones the anchor paragraph's formatting\n# (pStyle \"Listeavsnitt\", numPr ilvl 0 / numId 1, rPr lang en-GB), so the\n# new paragraph automatically becomes a list item in the same list.\n$anchor.Range.InsertParagraphAfter()\n$first = $anchor.Next()\n$first.Range.Text = \"(Confirm the performance of the current sensor at low currents.)\"\n\n$first.Range.InsertParagraphAfter()\n$second = $first.Next()\n$second.Range.Text = \"Testing the new constant voltage state\"\n"}
